$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("final_fail")
$ws2 = $wb.Worksheets.Item("final_gifted")

$ws1.Cells.Item(7,4).Value = $false
$ws1.Cells.Item(7,10).Value = 5
$ws1.Cells.Item(8,2).Value = $false
$ws1.Cells.Item(8,10).Value = 5
$ws1.Cells.Item(9,4).Value = $false
$ws1.Cells.Item(9,6).Value = $true
$ws1.Cells.Item(10,4).Value = $false
$ws1.Cells.Item(10,5).Value = $true
$ws1.Cells.Item(11,4).Value = $true
$ws1.Cells.Item(11,5).Value = $false
$ws1.Cells.Item(12,4).Value = $true
$ws1.Cells.Item(12,6).Value = $false
$ws1.Cells.Item(18,4).Value = $true
$ws1.Cells.Item(18,10).Value = 5
$ws1.Cells.Item(19,5).Value = $true
$ws1.Cells.Item(19,6).Value = $false
$ws1.Cells.Item(20,2).Value = $false
$ws1.Cells.Item(20,6).Value = $true
$ws1.Cells.Item(21,2).Value = $true
$ws1.Cells.Item(21,3).Value = $true
$ws1.Cells.Item(21,4).Value = $true
$ws1.Cells.Item(21,5).Value = $false
$ws1.Cells.Item(21,6).Value = $false
$ws1.Cells.Item(21,10).Value = 4
$ws1.Cells.Item(22,2).Value = $true
$ws1.Cells.Item(22,4).Value = $false
$ws1.Cells.Item(22,6).Value = $true
$ws1.Cells.Item(22,10).Value = 4
$ws1.Cells.Item(27,3).Value = $true
$ws1.Cells.Item(27,4).Value = $false
$ws1.Cells.Item(28,3).Value = $true
$ws1.Cells.Item(28,4).Value = $false
$ws1.Cells.Item(29,3).Value = $false
$ws1.Cells.Item(29,4).Value = $true
$ws1.Cells.Item(30,3).Value = $false
$ws1.Cells.Item(30,4).Value = $true
$ws1.Cells.Item(32,3).Value = $true
$ws1.Cells.Item(32,10).Value = 2
$ws1.Cells.Item(33,3).Value = $true
$ws1.Cells.Item(33,10).Value = 2
$ws1.Cells.Item(34,3).Value = $true
$ws1.Cells.Item(34,10).Value = 2
$ws1.Cells.Item(35,3).Value = $true
$ws1.Cells.Item(35,10).Value = 2
$ws2.Cells.Item(2,1).Value = "On/off campus click ratio"
$ws2.Cells.Item(3,1).Value = "Total time online (min)"
$ws2.Cells.Item(5,1).Value = "Clicks per day"
$ws2.Cells.Item(6,1).Value = "Clicks (% of course total)"
$ws2.Cells.Item(7,1).Value = "Average grade of assignments"
$ws2.Cells.Item(8,1).Value = "Resources viewed"
$ws2.Cells.Item(22,1).Value = "Clicks on folder"
$ws2.Cells.Item(23,1).Value = "Clicks on forum"
$ws2.Cells.Item(36,1).Value = "Submissions (% of course total)"
$ws2.Cells.Item(37,1).Value = "Start of Session 6 (%)"
$ws2.Cells.Item(9,1).Value = "Days with no interaction"
$ws2.Cells.Item(10,1).Value = "Clicks per session"
$ws2.Cells.Item(11,1).Value = "Days with no interaction (%)"
$ws2.Cells.Item(12,1).Value = "Clicks on course"
$ws2.Cells.Item(13,1).Value = "Start of Session 3 (%)"
$ws2.Cells.Item(14,1).Value = "Start of Session 2 (%)"
$ws2.Cells.Item(14,3).Value = $true
$ws2.Cells.Item(14,10).Value = 5
$ws2.Cells.Item(15,1).Value = "Start of Session 1 (%)"
$ws2.Cells.Item(15,3).Value = $true
$ws2.Cells.Item(15,10).Value = 5
$ws2.Cells.Item(16,1).Value = "Largest period of inactivity (h)"
$ws2.Cells.Item(16,3).Value = $true
$ws2.Cells.Item(16,10).Value = 5
$ws2.Cells.Item(17,1).Value = "Start of Session 4 (%)"
$ws2.Cells.Item(17,3).Value = $true
$ws2.Cells.Item(17,5).Value = $false
$ws2.Cells.Item(18,1).Value = "Clicks on campus"
$ws2.Cells.Item(18,2).Value = $false
$ws2.Cells.Item(18,5).Value = $true
$ws2.Cells.Item(19,1).Value = "Assignments viewed"
$ws2.Cells.Item(19,2).Value = $true
$ws2.Cells.Item(19,3).Value = $true
$ws2.Cells.Item(19,5).Value = $false
$ws2.Cells.Item(19,6).Value = $false
$ws2.Cells.Item(20,1).Value = "Number of clicks"
$ws2.Cells.Item(20,4).Value = $false
$ws2.Cells.Item(20,5).Value = $true
$ws2.Cells.Item(20,6).Value = $true
$ws2.Cells.Item(20,10).Value = 3
$ws2.Cells.Item(21,1).Value = "Files downloaded"
$ws2.Cells.Item(21,2).Value = $false
$ws2.Cells.Item(21,4).Value = $true
$ws2.Cells.Item(22,1).Value = "Clicks on folder"
$ws2.Cells.Item(23,1).Value = "Clicks on forum"
$ws2.Cells.Item(24,1).Value = "Start of Session 7 (%)"
$ws2.Cells.Item(24,2).Value = $true
$ws2.Cells.Item(24,6).Value = $false
$ws2.Cells.Item(25,1).Value = "Start of Session 5 (%)"
$ws2.Cells.Item(25,4).Value = $false
$ws2.Cells.Item(25,6).Value = $true
$ws2.Cells.Item(26,1).Value = "Links viewed"
$ws2.Cells.Item(26,4).Value = $true
$ws2.Cells.Item(26,10).Value = 2
$ws2.Cells.Item(27,1).Value = "Assignments submitted"
$ws2.Cells.Item(28,1).Value = "Number of sessions"
$ws2.Cells.Item(29,1).Value = "Discussions viewed"
$ws2.Cells.Item(30,1).Value = "Quizzes started"
$ws2.Cells.Item(31,1).Value = "Forum posts"
$ws2.Cells.Item(32,1).Value = "Start of Session 10 (%)"
$ws2.Cells.Item(33,1).Value = "Number of days"
$ws2.Cells.Item(34,1).Value = "Start of Session 9 (%)"
$ws2.Cells.Item(35,1).Value = "Start of Session 8 (%)"
$ws2.Cells.Item(36,1).Value = "Submissions (% of course total)"
$ws2.Cells.Item(37,1).Value = "Start of Session 6 (%)"
